$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Re-format the date column (Investment Date) cells: drop the
# "Normal 4" cell-style link and switch to the plain default font while
# keeping the date number format. ---
$ws.Range("C2").Style = "Normal"
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# Amount / Quantity / Instrument cells on the existing rows lose their
# explicit formatting and fall back to the plain default style.
$ws.Range("D2:E3").Style = "Normal"
$ws.Range("H2:H3").Style = "Normal"
$excel.CutCopyMode = $false

# --- Add the new data row (row 4) ---
$ws.Range("A4").Value = "Fund X"
$ws.Range("A4").Style = "Normal"

$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Ego Pvt Ltd"

$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value2 = 44880

$ws.Range("D4").Value = 300000
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = 150000
$ws.Range("E4").Style = "Normal"

$ws.Range("F4").Value = "Test"

$ws.Range("H4").Value = "Equity"
$ws.Range("H4").Style = "Normal"

$ws.Range("I4").Value = "INR"

$excel.CutCopyMode = $false
$ws.Range("E4").Select()
